# Iraq League workbook update - 23-02-2024
# This script:
#  1) Fixes 12 pairs of adjacent rows whose match data (everything except
#     the running index in col A, Div/Div Original Name in C/D, and the
#     Date in col E) had been attached to the wrong row - swap col B and
#     the F:AC block between each pair.
#  2) Fills in the result (FTHG/FTAG/FTR) and closing odds for the match
#     in row 207, which had kicked off since the data was last pulled, and
#     corrects its Date.
#  3) Appends 3 new fixtures/results (rows 208-210) that were added to the
#     source feed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Swap mis-attached match data between adjacent row pairs ---------
$pairs = @(
  @(16,17),
  @(25,26),
  @(39,40),
  @(64,65),
  @(80,81),
  @(87,88),
  @(137,138),
  @(141,142),
  @(152,153),
  @(156,157),
  @(161,162),
  @(186,187)
)

foreach ($p in $pairs) {
  $r1 = $p[0]
  $r2 = $p[1]

  $b1 = $ws.Range("B$r1").Value()
  $b2 = $ws.Range("B$r2").Value()
  $ws.Range("B$r1").Value = $b2
  $ws.Range("B$r2").Value = $b1

  $row1 = $ws.Range("F${r1}:AC${r1}").Value()
  $row2 = $ws.Range("F${r2}:AC${r2}").Value()
  $ws.Range("F${r1}:AC${r1}").Value = $row2
  $ws.Range("F${r2}:AC${r2}").Value = $row1
}

# --- 2) Row 207: match has been played, fill in result + closing odds ---
$ws.Range("E207").Value = 45344.33333333334
$ws.Range("H207").Value = 0
$ws.Range("I207").Value = 2
$ws.Range("J207").Value = "A"
$ws.Range("N207").Value = 3.25
$ws.Range("O207").Value = 3
$ws.Range("R207").Value = 1.9
$ws.Range("S207").Value = 1.9
$ws.Range("T207").Value = 2
$ws.Range("W207").Value = -1
$ws.Range("X207").Value = -1
$ws.Range("Y207").Value = 1.15
$ws.Range("Z207").Value = -1
$ws.Range("AA207").Value = 0.8999999999999999
$ws.Range("AB207").Value = 0
$ws.Range("AC207").Value = -0

# --- 3) New fixtures appended to the base ---------------------------------

# Row 208 - played, full stats
$ws.Range("A208").Value = 206
$ws.Range("B208").Value = 7864431
$ws.Range("C208").Value = "Iraq League"
$ws.Range("D208").Value = "Iraq League"
$ws.Range("E208").Value = 45344.4375
$ws.Range("F208").Value = "Al Najaf"
$ws.Range("G208").Value = "Erbil SC"
$ws.Range("H208").Value = 1
$ws.Range("I208").Value = 1
$ws.Range("J208").Value = "D"
$ws.Range("K208").Value = 2.2
$ws.Range("L208").Value = 2.9
$ws.Range("M208").Value = 3.2
$ws.Range("N208").Value = 2.25
$ws.Range("O208").Value = 2.9
$ws.Range("P208").Value = 3.1
$ws.Range("Q208").Value = -0.25
$ws.Range("R208").Value = 2
$ws.Range("S208").Value = 1.8
$ws.Range("T208").Value = 1.75
$ws.Range("U208").Value = 1.775
$ws.Range("V208").Value = 2.025
$ws.Range("W208").Value = -1
$ws.Range("X208").Value = 1.9
$ws.Range("Y208").Value = -1
$ws.Range("Z208").Value = -0.5
$ws.Range("AA208").Value = 0.4
$ws.Range("AB208").Value = 0.3875
$ws.Range("AC208").Value = -0.5

# Row 209 - not yet played (no FTHG/FTAG/FTR, no PL_AhOver/PL_AhUnder)
$ws.Range("A209").Value = 207
$ws.Range("B209").Value = 7870836
$ws.Range("C209").Value = "Iraq League"
$ws.Range("D209").Value = "Iraq League"
$ws.Range("E209").Value = 45345.33333333334
$ws.Range("F209").Value = "Al Hudod"
$ws.Range("G209").Value = "Naft Maysan"
$ws.Range("K209").Value = 2.6
$ws.Range("L209").Value = 2.7
$ws.Range("M209").Value = 2.8
$ws.Range("N209").Value = 2.45
$ws.Range("O209").Value = 2.6
$ws.Range("P209").Value = 3.1
$ws.Range("Q209").Value = -0.25
$ws.Range("R209").Value = 2.1
$ws.Range("S209").Value = 1.7
$ws.Range("T209").Value = 1.75
$ws.Range("U209").Value = 1.825
$ws.Range("V209").Value = 1.975
$ws.Range("W209").Value = 0
$ws.Range("X209").Value = 0
$ws.Range("Y209").Value = 0
$ws.Range("Z209").Value = 0
$ws.Range("AA209").Value = 0

# Row 210 - not yet played
$ws.Range("A210").Value = 208
$ws.Range("B210").Value = 7873977
$ws.Range("C210").Value = "Iraq League"
$ws.Range("D210").Value = "Iraq League"
$ws.Range("E210").Value = 45345.4375
$ws.Range("F210").Value = "Duhok"
$ws.Range("G210").Value = "Karbalaa FC"
$ws.Range("K210").Value = 1.8
$ws.Range("L210").Value = 3
$ws.Range("M210").Value = 4.5
$ws.Range("N210").Value = 1.8
$ws.Range("O210").Value = 3
$ws.Range("P210").Value = 4.5
$ws.Range("Q210").Value = -0.5
$ws.Range("R210").Value = 1.825
$ws.Range("S210").Value = 1.975
$ws.Range("T210").Value = 2
$ws.Range("U210").Value = 2
$ws.Range("V210").Value = 1.8
$ws.Range("W210").Value = 0
$ws.Range("X210").Value = 0
$ws.Range("Y210").Value = 0
$ws.Range("Z210").Value = 0
$ws.Range("AA210").Value = 0
